# "Generate Report for Handback"
#
# The 4c2112cc-...md file has now been handed back (status flips from
# "Ready for handoff" to "Handed back: in sync with en-US"), so the
# generated report now lists it first (row 2) on every sheet, pushing
# 455823c7-...md (still "Ready for handoff") down to row 3. The zh-cn /
# de-de detail sheets also gain freshly-populated "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns for the
# file that was just handed back.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(2,1).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$ws.Cells.Item(2,2).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "2016-42-13 20:42:27"

$ws.Cells.Item(3,1).Value = "455823c7-66b0-462f-8e16-75152ac7a3f5.md"
$ws.Cells.Item(3,2).Value = "Ready for handoff"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "2016-41-13 20:41:33"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fb3e0eae4f67802deade0be3256c97e2d7aabe6f/e2e/455823c7-66b0-462f-8e16-75152ac7a3f5.md", $null, $null, "455823c7-66b0-462f-8e16-75152ac7a3f5.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2: 4c2112cc-... file, now handed back
$ws.Cells.Item(2,1).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-13 20:42:24"
$ws.Cells.Item(2,6).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$ws.Cells.Item(2,7).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf"
$ws.Cells.Item(2,8).Value = "2016-03-13 20:45:35"
$ws.Cells.Item(2,9).Value = "Include"

# Row 3: 455823c7-... file, still ready for handoff
$ws.Cells.Item(3,1).Value = "455823c7-66b0-462f-8e16-75152ac7a3f5.md"
$ws.Cells.Item(3,2).Value = ".md"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.zh-cn.xlf"
$ws.Cells.Item(3,5).Value = "2016-03-13 20:39:17"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75ea7a263ec19184071310fd35e9f0fcbccc2f9b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75ea7a263ec19184071310fd35e9f0fcbccc2f9b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fb3e0eae4f67802deade0be3256c97e2d7aabe6f/e2e/455823c7-66b0-462f-8e16-75152ac7a3f5.md", $null, $null, "455823c7-66b0-462f-8e16-75152ac7a3f5.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fb3e0eae4f67802deade0be3256c97e2d7aabe6f/e2e/455823c7-66b0-462f-8e16-75152ac7a3f5.md", $null, $null, ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/058bb91b06be86dc76dd8ec7eeee514dbe19b691/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.zh-cn.xlf", $null, $null, "455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2: 4c2112cc-... file, now handed back
$ws.Cells.Item(2,1).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$ws.Cells.Item(2,2).Value = ".md"
$ws.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws.Cells.Item(2,4).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf"
$ws.Cells.Item(2,5).Value = "2016-03-13 20:42:27"
$ws.Cells.Item(2,6).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md"
$ws.Cells.Item(2,7).Value = "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf"
$ws.Cells.Item(2,8).Value = "2016-03-13 20:45:42"
$ws.Cells.Item(2,9).Value = "Include"

# Row 3: 455823c7-... file, still ready for handoff
$ws.Cells.Item(3,1).Value = "455823c7-66b0-462f-8e16-75152ac7a3f5.md"
$ws.Cells.Item(3,2).Value = ".md"
$ws.Cells.Item(3,3).Value = "Ready for handoff"
$ws.Cells.Item(3,4).Value = "455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.de-de.xlf"
$ws.Cells.Item(3,5).Value = "2016-03-13 20:41:33"
$ws.Cells.Item(3,8).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(3,9).Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c160fe6f3bb3c395655aec1d160ba1aaa459f1fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/7591127eb3c2472b8a55019b461519b26debee80/e2e/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.md")
$ws.Hyperlinks.Add($ws.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c160fe6f3bb3c395655aec1d160ba1aaa459f1fa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf", $null, $null, "4c2112cc-b53a-4e04-9d2c-d5b1abfd772f.10de492e28aa8324ca31d2fcf443f901ddfb7084.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/fb3e0eae4f67802deade0be3256c97e2d7aabe6f/e2e/455823c7-66b0-462f-8e16-75152ac7a3f5.md", $null, $null, "455823c7-66b0-462f-8e16-75152ac7a3f5.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/fb3e0eae4f67802deade0be3256c97e2d7aabe6f/e2e/455823c7-66b0-462f-8e16-75152ac7a3f5.md", $null, $null, ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/49abf0c649a365703bc551d61f9fa0ee2d0d16c2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.de-de.xlf", $null, $null, "455823c7-66b0-462f-8e16-75152ac7a3f5.9ba1e1f41680c90f21fc850b354948b9281e557e.de-de.xlf")
